$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new Kaspa buy recorded on 2025-07-04 as row 38.
# Leading apostrophe forces the date-looking string to be stored as text
# (matching the existing "MM/DD/YYYY" text cells used elsewhere in column A),
# then ClearFormats drops the resulting quote-prefix style so the cell keeps
# the sheet's default (unstyled) formatting.
$ws.Range("A38").Value = "'07/04/2025"
$ws.Range("A38").ClearFormats()

$ws.Range("B38").Value = 629.2050000000017
$ws.Range("C38").Value = 0.0794653570775818
$ws.Range("D38").Value = 50
